$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.828.15'
$ws.Range('E2').Value = '  -1.59%  '
$ws.Range('D3').Value = '1.874.37'
$ws.Range('E3').Value = '  -1.84%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = "'301.45"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.09%  '
$ws.Range('D6').Value = "'1.002"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.06%  '
$ws.Range('D7').Value = "'0.5350"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.97%  '
$ws.Range('D8').Value = "'0.3747"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.04%  '
$ws.Range('D9').Value = "'0.07195"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.60%  '
$ws.Range('D10').Value = "'21.57"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.29%  '
$ws.Range('D11').Value = "'0.8901"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.88%  '
$ws.Range('D12').Value = "'0.08193"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.23%  '
$ws.Range('D13').Value = '1.874.43'
$ws.Range('E13').Value = '  +6.97%  '
$ws.Range('D14').Value = "'93.27"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.14%  '
$ws.Range('D15').Value = "'5.318"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.03%  '
$ws.Range('D16').Value = "'1.002"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.17%  '
$ws.Range('D17').Value = "'14.83"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.37%  '
$ws.Range('D18').Value = "'0.000008534"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.79%  '
$ws.Range('D19').Value = "'1.002"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.02%  '
$ws.Range('D20').Value = '26.867.08'
$ws.Range('E20').Value = '  -1.56%  '
$ws.Range('E21').Value = '  -2.73%  '
$ws.Range('D22').Value = "'10.61"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.84%  '
$ws.Range('D23').Value = "'6.394"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.67%  '
$ws.Range('D24').Value = "'2.289"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.48%  '
$ws.Range('D25').Value = "'146.22"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.67%  '
$ws.Range('D26').Value = "'18.10"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.84%  '
$ws.Range('D27').Value = "'1.731"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.57%  '
$ws.Range('D28').Value = "'114.05"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.58%  '
$ws.Range('D29').Value = "'4.715"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.96%  '
$ws.Range('D30').Value = "'4.614"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.49%  '
$ws.Range('D31').Value = "'0.09116"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.27%  '
$ws.Range('D32').Value = "'0.8076"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.24%  '
$ws.Range('D33').Value = "'0.05015"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.34%  '
$ws.Range('D34').Value = "'1.174"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.00%  '
$ws.Range('D35').Value = "'2.957"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.03%  '
$ws.Range('D36').Value = "'0.6062"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.35%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').Value = "'2.643"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.72%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').Value = "'3.212"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.81%  '
$ws.Range('D39').Value = "'0.01956"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.45%  '
$ws.Range('D40').Value = "'1.072"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.06%  '
$ws.Range('D41').Value = "'6.599"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.40%  '
$ws.Range('D42').Value = "'8.875"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.03%  '
$ws.Range('D43').Value = "'0.5129"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.79%  '
$ws.Range('D44').Value = "'114.95"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.08%  '
$ws.Range('D45').Value = "'0.1497"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.72%  '
$ws.Range('D46').Value = "'1.002"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.05%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = "'10.00"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.53%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').Value = "'1.642"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.17%  '
$ws.Range('D49').Value = "'37.59"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.93%  '
$ws.Range('D50').Value = "'0.06083"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.50%  '
$ws.Range('E51').Value = '  -3.07%  '
